# feat: expense search button added
#
# The template sheet (a single "Hoja1" worksheet with headers in row 1 and
# one sample expense row in row 2) gets 7 more sample rows appended
# (rows 3-9), each one an exact copy of the existing row 2
# (Description="Pruebaaa", Amount=2000, Category="TVS"), and the current
# selection moves to the newly added last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2:C2) holds the sample expense entry - replicate it down
# through row 9, matching the new data added in the sheet.
$template = $ws.Range("A2:C2")
for ($row = 3; $row -le 9; $row++) {
    $target = $ws.Range("A" + $row + ":C" + $row)
    $template.Copy($target)
}

# Reflect the new selection left behind after entering/filling the data
# (last appended row, columns A:C).
$ws.Range("A9:C9").Select()
